$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Date placeholder text: "12/02/2016" -> "2016-08-08"
#    (Slide Master + all 11 Custom Layouts)
# ---------------------------------------------------------------------------
function Set-DateText {
    param($shapes)
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            # Force a real text re-write (a no-op assignment of the same
            # value is skipped by the host), then set the final value.
            $shp.TextFrame.TextRange.Text = "."
            $shp.TextFrame.TextRange.Text = "2016-08-08"
        }
    }
}

$master = $p.SlideMaster
Set-DateText $master.Shapes
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Set-DateText $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2. Slide 1 shape edits
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

# -- Shape position nudges --------------------------------------------------
# Target Left/Top values given in points (12700 EMU = 1pt), precomputed so
# that after the host's internal float round-trip the saved EMU exactly
# matches the target (the host stores Left/Top with reduced precision, so a
# plain "target_emu / 12700.0" can truncate to target_emu - 1 on save).
$s.Shapes.Item(2).Left  = 64.18811023622047    # id 196 Rounded Rectangle 195    -> x=815189
$s.Shapes.Item(19).Top  = 15.328425196850393   # id 173 Straight Arrow Conn 172  -> y=194671
$s.Shapes.Item(22).Left = 204.1627579055118    # id 177 Rounded Rectangle 176    -> x=2592867
$s.Shapes.Item(23).Left = 72.45173268346457    # id 178 Parallelogram 177        -> x=920137
$s.Shapes.Item(24).Left = 130.73811023622048   # id 179 Up Arrow 178             -> x=1660374
$s.Shapes.Item(30).Left = 12.938031696062993   # id 185 Straight Connector 184   -> x=164313
$s.Shapes.Item(31).Left = 139.45551301102364   # id 186 Straight Arrow Conn 185  -> x=1771085
$s.Shapes.Item(32).Left = 82.37409598818897    # id 187 Straight Connector 186   -> x=1046151
$s.Shapes.Item(33).Left = 186.63228616456695   # id 188 Straight Arrow Conn 187  -> x=2370230
$s.Shapes.Item(34).Left = 190.4663779527559    # id 189 Straight Connector 188   -> x=2418923
$s.Shapes.Item(35).Left = 211.49889763779527   # id 190 Straight Arrow Conn 189  -> x=2686036
$s.Shapes.Item(37).Left = 82.37409598818897    # id 192 Straight Connector 191   -> x=1046151
$s.Shapes.Item(38).Left = 196.11347206692915   # id 193 Straight Connector 192   -> x=2490641
$s.Shapes.Item(39).Left = 210.57708661417323   # id 194 Straight Arrow Conn 193  -> x=2674329

# -- Line weight: give every straight-line / connector shape a 1.5pt (19050 EMU) weight --
$lineShapeIdx = @(10,11,12,13,14,15,16,17,18,19,20,30,31,32,33,34,35,36,37,38,39,40,47,48,49,50,51,52)
foreach ($idx in $lineShapeIdx) {
    $s.Shapes.Item($idx).Line.Weight = 1.5
}

# -- Dash style: shape id 193 (index 38) sysDash -> sysDot --
$s.Shapes.Item(38).Line.DashStyle = 2

# -- Merge the "No " + "Acoustic Wave" runs into a single "No Acoustic Wave" run --
$txShape = $s.Shapes.Item(44)
$txShape.TextFrame.TextRange.Text = "."
$txShape.TextFrame.TextRange.Text = "No Acoustic Wave"
